$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current header row (row 2) to make room
# for the new "order date" / "collect date" lines. This pushes the
# existing header + data rows down by two (old row 2 -> row 4, etc.) and
# carries their number formatting/merges along automatically.
$ws.Rows.Item(2).Resize(2).Insert()

# Fill in the new label/value pairs.
$ws.Range("A2").Value = "Data zamówienia:"
$ws.Range("B2").Value = "22.4.2020"
$ws.Range("A3").Value = "Data odbioru:"
$ws.Range("B3").Value = "30.4.2020"

# Give the two new rows a uniform look (same base formatting for the
# label and the value cell) by pulling the format from the title cell
# A1, then shrinking the font down to 12pt / non-bold for these rows.
$labelRange = $ws.Range("A2:C3")
$ws.Range("A1").Copy()
$labelRange.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
$labelRange.Font.Size = 12
$labelRange.Font.Bold = $false

# Merge the value cell across B:C, same as the other label/value pairs
# lower in the sheet.
$ws.Range("B2:C2").Merge()
$ws.Range("B3:C3").Merge()
